# Update "想去人数" (F column) counts per worksheet, matching the scraped
# data refresh captured in the commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 312
$ws1.Range("F7").Value  = 1170
$ws1.Range("F8").Value  = 443
$ws1.Range("F9").Value  = 7104
$ws1.Range("F13").Value = 7971
$ws1.Range("F15").Value = 51
$ws1.Range("F16").Value = 5506
$ws1.Range("F18").Value = 2403
$ws1.Range("F23").Value = 78
$ws1.Range("F25").Value = 375
$ws1.Range("F26").Value = 255
$ws1.Range("F28").Value = 2342
$ws1.Range("F31").Value = 81
$ws1.Range("F32").Value = 143
$ws1.Range("F33").Value = 580
$ws1.Range("F34").Value = 5
$ws1.Range("F39").Value = 2309
$ws1.Range("F42").Value = 6

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 92
$ws2.Range("F3").Value = 76
$ws2.Range("F4").Value = 62
$ws2.Range("F5").Value = 3

# Sheet "全部类型" (All types, merged listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 312
$ws4.Range("F7").Value  = 92
$ws4.Range("F9").Value  = 1170
$ws4.Range("F10").Value = 443
$ws4.Range("F11").Value = 7105
$ws4.Range("F15").Value = 7971
$ws4.Range("F17").Value = 51
$ws4.Range("F18").Value = 5506
$ws4.Range("F20").Value = 2403
$ws4.Range("F24").Value = 78
$ws4.Range("F25").Value = 76
$ws4.Range("F27").Value = 62
$ws4.Range("F28").Value = 375
$ws4.Range("F30").Value = 2342
$ws4.Range("F33").Value = 81
$ws4.Range("F34").Value = 143
$ws4.Range("F35").Value = 3
$ws4.Range("F36").Value = 580
$ws4.Range("F37").Value = 5
$ws4.Range("F43").Value = 2309
$ws4.Range("F47").Value = 6
